# Fix bug with db writing. Add 'working time' feature in db entity
# -> append a new contact row to the phone base sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 18: name in column A, phone number (kept as text) in column B.
$ws.Range("A18").Value = "Маша Тестировщик"
$ws.Range("B18").Value = "+380954121725"

# Move/restore the active selection to A19, just below the newly added row,
# matching where Excel leaves the cursor after the data entry.
$ws.Range("A19").Select()
